$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-5 with new measurement values ---
$ws.Cells.Item(2, 1).Value = 45051.50694444445
$ws.Cells.Item(2, 2).Value = 15.854
$ws.Cells.Item(2, 3).Value = 10.494
$ws.Cells.Item(2, 4).Value = 3.735
$ws.Cells.Item(2, 5).Value = 33.762
$ws.Cells.Item(2, 6).Value = 26.097
$ws.Cells.Item(2, 7).Value = 12.381
$ws.Cells.Item(2, 8).Value = 37.784
$ws.Cells.Item(2, 9).Value = 19.197
$ws.Cells.Item(2, 10).Value = 7.826
$ws.Cells.Item(2, 11).Value = 11.635
$ws.Cells.Item(2, 12).Value = 13.329
$ws.Cells.Item(2, 13).Value = 13.955
$ws.Cells.Item(2, 14).Value = 3.981
$ws.Cells.Item(2, 15).Value = 12.407
$ws.Cells.Item(2, 16).Value = 17.155
$ws.Cells.Item(2, 17).Value = 10.858
$ws.Cells.Item(2, 18).Value = 3.133
$ws.Cells.Item(2, 19).Value = 2.01
$ws.Cells.Item(2, 20).Value = 181.475
$ws.Cells.Item(2, 21).Value = 34.412
$ws.Cells.Item(2, 22).Value = 11.452
$ws.Cells.Item(2, 23).Value = 22.298
$ws.Cells.Item(2, 24).Value = 11.285
$ws.Cells.Item(2, 25).Value = 3.17
$ws.Cells.Item(2, 26).Value = 19.349
$ws.Cells.Item(2, 27).Value = 10.115
$ws.Cells.Item(2, 28).Value = 9.106
$ws.Cells.Item(2, 29).Value = 10.918
$ws.Cells.Item(2, 30).Value = 14.036
$ws.Cells.Item(2, 31).Value = 3.313
$ws.Cells.Item(2, 32).Value = 33.905
$ws.Cells.Item(2, 33).Value = 6.108
$ws.Cells.Item(2, 34).Value = 14.317
$ws.Cells.Item(3, 1).Value = 45051.51388888889
$ws.Cells.Item(3, 2).Value = 23.06
$ws.Cells.Item(3, 3).Value = 16.639
$ws.Cells.Item(3, 4).Value = 1.999
$ws.Cells.Item(3, 5).Value = 50.025
$ws.Cells.Item(3, 6).Value = 40.321
$ws.Cells.Item(3, 7).Value = 18.065
$ws.Cells.Item(3, 8).Value = 68.01000000000001
$ws.Cells.Item(3, 9).Value = 27.923
$ws.Cells.Item(3, 10).Value = 12.197
$ws.Cells.Item(3, 11).Value = 17.981
$ws.Cells.Item(3, 12).Value = 20.011
$ws.Cells.Item(3, 13).Value = 21.099
$ws.Cells.Item(3, 14).Value = 5.797
$ws.Cells.Item(3, 15).Value = 18.046
$ws.Cells.Item(3, 16).Value = 25.524
$ws.Cells.Item(3, 17).Value = 15.422
$ws.Cells.Item(3, 18).Value = 1.534
$ws.Cells.Item(3, 19).Value = 1.266
$ws.Cells.Item(3, 20).Value = 267.37
$ws.Cells.Item(3, 21).Value = 50.451
$ws.Cells.Item(3, 22).Value = 16.657
$ws.Cells.Item(3, 23).Value = 33.61
$ws.Cells.Item(3, 24).Value = 17.523
$ws.Cells.Item(3, 25).Value = 3.074
$ws.Cells.Item(3, 26).Value = 33.518
$ws.Cells.Item(3, 27).Value = 14.713
$ws.Cells.Item(3, 28).Value = 13.126
$ws.Cells.Item(3, 29).Value = 15.495
$ws.Cells.Item(3, 30).Value = 21.036
$ws.Cells.Item(3, 31).Value = 1.246
$ws.Cells.Item(3, 32).Value = 62.13
$ws.Cells.Item(3, 33).Value = 9.242000000000001
$ws.Cells.Item(3, 34).Value = 20.825
$ws.Cells.Item(4, 1).Value = 45051.52083333334
$ws.Cells.Item(4, 2).Value = 9.609
$ws.Cells.Item(4, 3).Value = 6.767
$ws.Cells.Item(4, 4).Value = 1.08
$ws.Cells.Item(4, 5).Value = 20.844
$ws.Cells.Item(4, 6).Value = 16.495
$ws.Cells.Item(4, 7).Value = 7.488
$ws.Cells.Item(4, 8).Value = 33.32
$ws.Cells.Item(4, 9).Value = 11.635
$ws.Cells.Item(4, 10).Value = 5.023
$ws.Cells.Item(4, 11).Value = 7.262
$ws.Cells.Item(4, 12).Value = 8.343
$ws.Cells.Item(4, 13).Value = 8.792999999999999
$ws.Cells.Item(4, 14).Value = 2.419
$ws.Cells.Item(4, 15).Value = 7.519
$ws.Cells.Item(4, 16).Value = 10.594
$ws.Cells.Item(4, 17).Value = 6.585
$ws.Cells.Item(4, 18).Value = 0.955
$ws.Cells.Item(4, 19).Value = 0.627
$ws.Cells.Item(4, 20).Value = 107.148
$ws.Cells.Item(4, 21).Value = 21.188
$ws.Cells.Item(4, 22).Value = 6.941
$ws.Cells.Item(4, 23).Value = 13.965
$ws.Cells.Item(4, 24).Value = 7.234
$ws.Cells.Item(4, 25).Value = 1.447
$ws.Cells.Item(4, 26).Value = 15.626
$ws.Cells.Item(4, 27).Value = 6.131
$ws.Cells.Item(4, 28).Value = 5.549
$ws.Cells.Item(4, 29).Value = 6.541
$ws.Cells.Item(4, 30).Value = 8.76
$ws.Cells.Item(4, 31).Value = 0.766
$ws.Cells.Item(4, 32).Value = 30.676
$ws.Cells.Item(4, 33).Value = 3.773
$ws.Cells.Item(4, 34).Value = 8.678000000000001
$ws.Cells.Item(5, 1).Value = 45051.52777777778
$ws.Cells.Item(5, 2).Value = 11.53
$ws.Cells.Item(5, 3).Value = 8.33
$ws.Cells.Item(5, 4).Value = 0.93
$ws.Cells.Item(5, 5).Value = 25.05
$ws.Cells.Item(5, 6).Value = 20.14
$ws.Cells.Item(5, 7).Value = 9.01
$ws.Cells.Item(5, 8).Value = 35.3
$ws.Cells.Item(5, 9).Value = 13.96
$ws.Cells.Item(5, 10).Value = 6.1
$ws.Cells.Item(5, 11).Value = 8.94
$ws.Cells.Item(5, 12).Value = 10.04
$ws.Cells.Item(5, 13).Value = 10.62
$ws.Cells.Item(5, 14).Value = 2.9
$ws.Cells.Item(5, 15).Value = 9.02
$ws.Cells.Item(5, 16).Value = 12.75
$ws.Cells.Item(5, 17).Value = 7.76
$ws.Cells.Item(5, 18).Value = 0.75
$ws.Cells.Item(5, 19).Value = 0.58
$ws.Cells.Item(5, 20).Value = 130.01
$ws.Cells.Item(5, 21).Value = 25.22
$ws.Cells.Item(5, 22).Value = 8.33
$ws.Cells.Item(5, 23).Value = 16.78
$ws.Cells.Item(5, 24).Value = 8.779999999999999
$ws.Cells.Item(5, 25).Value = 1.52
$ws.Cells.Item(5, 26).Value = 16.98
$ws.Cells.Item(5, 27).Value = 7.36
$ws.Cells.Item(5, 28).Value = 6.59
$ws.Cells.Item(5, 29).Value = 7.76
$ws.Cells.Item(5, 30).Value = 10.55
$ws.Cells.Item(5, 31).Value = 0.55
$ws.Cells.Item(5, 32).Value = 32.07
$ws.Cells.Item(5, 33).Value = 4.61
$ws.Cells.Item(5, 34).Value = 10.41

# --- Remove the now-unused last data row ---
$ws.Rows.Item(6).Delete()

# --- Widen a subset of the per-junction columns from 7 to 8 characters ---
# (ColumnWidth 7.17 round-trips to a stored column width of exactly 8)
$wideCols = @(3,7,10,11,12,13,15,17,22,24,27,28,29,30,34)
foreach ($col in $wideCols) {
    $ws.Columns.Item($col).ColumnWidth = 7.17
}
